# Apply cryptos list price/volume updates (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.649.99"
$ws.Range("E2").Value = "  +3.72%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.418.65"
$ws.Range("E3").Value = "  +1.99%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.55"
$ws.Range("E5").Value = "  +4.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.16"
$ws.Range("E6").Value = "  +5.81%  "
$ws.Range("E7").Value = "  +2.02%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.527"
$ws.Range("E9").Value = "  +9.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.24"
$ws.Range("E10").Value = "  +2.56%  "
$ws.Range("E11").Value = "  +1.46%  "
$ws.Range("E12").Value = "  -1.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.59"
$ws.Range("E13").Value = "  +1.22%  "
$ws.Range("E14").Value = "  +2.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.796.97"
$ws.Range("E15").Value = "  +2.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.421.77"
$ws.Range("E16").Value = "  +2.27%  "
$ws.Range("E17").Value = "  +3.92%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "44.496.75"
$ws.Range("E18").Value = "  +3.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.19"
$ws.Range("E19").Value = "  +1.69%  "
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("E21").Value = "  +3.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.54"
$ws.Range("E22").Value = "  +0.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "241.71"
$ws.Range("E23").Value = "  +2.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.26"
$ws.Range("E24").Value = "  +3.21%  "
$ws.Range("E25").Value = "  +2.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.21"
$ws.Range("E27").Value = "  +3.03%  "
$ws.Range("E28").Value = "  -3.49%  "
$ws.Range("E29").Value = "  +1.90%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.29"
$ws.Range("E30").Value = "  +2.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "48.25"
$ws.Range("E31").Value = "  +0.88%  "
$ws.Range("E32").Value = "  +16.66%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.53"
$ws.Range("E33").Value = "  +10.84%  "
$ws.Range("E34").Value = "  +2.66%  "
$ws.Range("E35").Value = "  +0.24%  "
$ws.Range("E36").Value = "  +4.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.87"
$ws.Range("E37").Value = "  +2.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.43"
$ws.Range("E38").Value = "  +2.07%  "
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "125.65"
$ws.Range("E39").Value = "  -3.76%  "
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.84"
$ws.Range("E40").Value = "  -0.35%  "
$ws.Range("E41").Value = "  +1.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.18"
$ws.Range("E42").Value = "  -3.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.89"
$ws.Range("E43").Value = "  -0.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0288"
$ws.Range("E44").Value = "  +3.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.940.10"
$ws.Range("E45").Value = "  +0.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.12"
$ws.Range("E46").Value = "  -1.10%  "
$ws.Range("E47").Value = "  +7.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.13"
$ws.Range("E48").Value = "  -0.88%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.74"
$ws.Range("E49").Value = "  +15.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "75.00"
$ws.Range("E50").Value = "  +5.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.65"
$ws.Range("E51").Value = "  +4.85%  "
